# Apply the changes described in the commit "updated 4.0 files and mdl"
# to the "Fuel Prod Imp Exp Balancing Priorities" workbook.

$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

# --- About sheet: bump the "last updated" date shown in C1 -----------------
# 45294 (2024-01-03) -> 45379 (2024-03-28)
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: update the balancing-priority ranking for "hard coal" ---
# (row 3 : production / imports / exports priority order changed)
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# --- Restore FPIEBP as the active sheet and update the selected cell -------
$wsFPIEBP.Activate()
$wsFPIEBP.Range("E3").Select()
